# Update column G ("K") values on the active sheet for rows 2-26
# to reflect regenerated save_data (K instead of Strike#, regen std/mean,
# calc and write s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 0
    4  = 4
    5  = 4
    6  = 2
    7  = 3
    8  = 1
    9  = 3
    10 = 4
    11 = 7
    12 = 5
    13 = 4
    14 = 9
    15 = 3
    16 = 4
    17 = 5
    18 = 5
    19 = 7
    20 = 3
    21 = 1
    22 = 3
    23 = 1
    24 = 5
    25 = 2
    26 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
